$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Ccl12"
$ws.Range("C2").Value = "Ackr2"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 1.005755333333333
$ws.Range("H2").Value = 3.017266
$ws.Range("I2").Value = 0.01048729000197281
$ws.Range("J2").Value = 0.01048729000197281
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.26202
$ws.Range("N2").Value = 0.78606
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.26352801244
$ws.Range("R2").Value = 2.37175211196
$ws.Range("S2").Value = 0.01048729000197281
$ws.Range("T2").Value = 0.01048729000197281

# Row 3
$ws.Range("A3").Value = "Inflammatory-Mac"
$ws.Range("B3").Value = "Ccl12"
$ws.Range("C3").Value = "Ackr2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 38.755371
$ws.Range("H3").Value = 116.266113
$ws.Range("I3").Value = 0.4041130097356814
$ws.Range("J3").Value = 0.4041130097356814
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.26202
$ws.Range("N3").Value = 0.78606
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 10.15468230942
$ws.Range("R3").Value = 91.39214078478
$ws.Range("S3").Value = 0.4041130097356814
$ws.Range("T3").Value = 0.4041130097356814

# Row 4
$ws.Range("A4").Value = "Neutrophils"
$ws.Range("B4").Value = "Ccl12"
$ws.Range("C4").Value = "Ackr2"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 7.988471999999999
$ws.Range("H4").Value = 23.965416
$ws.Range("I4").Value = 0.08329801469605898
$ws.Range("J4").Value = 0.08329801469605898
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.26202
$ws.Range("N4").Value = 0.78606
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 2.093139433439999
$ws.Range("R4").Value = 18.83825490096
$ws.Range("S4").Value = 0.08329801469605898
$ws.Range("T4").Value = 0.08329801469605898

# Row 5
$ws.Range("A5").Value = "Resolving-Mac"
$ws.Range("B5").Value = "Ccl12"
$ws.Range("C5").Value = "Ackr2"
$ws.Range("D5").Value = "FAPs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 48.15271133333334
$ws.Range("H5").Value = 144.458134
$ws.Range("I5").Value = 0.5021016855662869
$ws.Range("J5").Value = 0.5021016855662868
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.26202
$ws.Range("N5").Value = 0.78606
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 1
$ws.Range("Q5").Value = 12.61697342356
$ws.Range("R5").Value = 113.55276081204
$ws.Range("S5").Value = 0.5021016855662869
$ws.Range("T5").Value = 0.5021016855662868
